# Weekly refresh of Fruta/Hortaliza data: the rows of this sheet keep the
# same set of (Fecha, Volumen, Precio mínimo, Precio máximo, Precio promedio
# ponderado, Precio $/Kg) records, but those values are reshuffled across the
# different market-day rows. Columns A,B,C,E-L,Q,R,T are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row -> source row (read the source row's values *before*
# any row has been overwritten, then apply them to the destination row).
$rowMap = @{}
$rowMap[2]  = 21
$rowMap[3]  = 23
$rowMap[4]  = 10
$rowMap[5]  = 33
$rowMap[6]  = 25
$rowMap[7]  = 27
$rowMap[8]  = 13
$rowMap[9]  = 16
$rowMap[10] = 22
$rowMap[11] = 7
$rowMap[12] = 15
$rowMap[13] = 14
$rowMap[14] = 32
$rowMap[15] = 11
$rowMap[16] = 29
$rowMap[17] = 3
$rowMap[18] = 9
$rowMap[19] = 8
$rowMap[20] = 30
$rowMap[21] = 6
$rowMap[22] = 12
$rowMap[23] = 2
$rowMap[24] = 18
$rowMap[25] = 24
$rowMap[26] = 31
$rowMap[27] = 4
$rowMap[28] = 5
$rowMap[29] = 17
$rowMap[30] = 34
$rowMap[31] = 20
$rowMap[32] = 26
$rowMap[33] = 28
$rowMap[34] = 19

# Columns whose values travel together with a row's record.
$cols = @(4, 13, 14, 15, 16, 19)

# Snapshot every source row's values first so that overwriting a row that is
# also used as someone else's source does not corrupt later reads.
$snapshot = @{}
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $vals = @{}
        foreach ($col in $cols) {
            $vals[$col] = $ws.Cells.Item($srcRow, $col).Value2
        }
        $snapshot[$srcRow] = $vals
    }
}

# Now write the snapshot values into their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value2 = $vals[$col]
    }
}
